$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, border, centered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record data for every data row (2-48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 98
    $ws.Cells.Item($r, 31).Value = 64
    $ws.Cells.Item($r, 32).Value = 0
}
